$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.302731113197012
$ws.Range("C2").Value = 0.1348401523616332
$ws.Range("D2").Value = 0.4636836715194335
$ws.Range("E2").Value = 0.1184377548800093
$ws.Range("G2").Value = 0.002784549938200505
$ws.Range("J2").Value = 0.02575336588086152
$ws.Range("K2").Value = 3.595999680862008
$ws.Range("L2").Value = 0.6933962085670231
$ws.Range("N2").Value = 7.390207188599476

$ws.Range("B3").Value = 4.261376840821526
$ws.Range("C3").Value = 0.1290371097969825
$ws.Range("D3").Value = 0.4624233527037518
$ws.Range("E3").Value = 0.1186677422826392
$ws.Range("G3").Value = 0.002790667944192828
$ws.Range("J3").Value = 0.02572025073381745
$ws.Range("K3").Value = 3.551250404210549
$ws.Range("L3").Value = 0.6905708056876421
$ws.Range("N3").Value = 7.315733859723991

$ws.Range("B4").Value = 4.238149426166899
$ws.Range("C4").Value = 0.1255620464430649
$ws.Range("D4").Value = 0.4618468495392847
$ws.Range("E4").Value = 0.1188476709181714
$ws.Range("G4").Value = 0.002794620825851842
$ws.Range("J4").Value = 0.02569982516832248
$ws.Range("K4").Value = 3.525674296202141
$ws.Range("L4").Value = 0.6891468745419189
$ws.Range("N4").Value = 7.270470312392803

$ws.Range("B5").Value = 4.229227594964982
$ws.Range("C5").Value = 0.1241679008563636
$ws.Range("D5").Value = 0.4616615387711107
$ws.Range("E5").Value = 0.1189307406176461
$ws.Range("G5").Value = 0.002796281217964777
$ws.Range("J5").Value = 0.02569147685076345
$ws.Range("K5").Value = 3.51572874336793
$ws.Range("L5").Value = 0.6886447434716985
$ws.Range("N5").Value = 7.252140536431597

$ws.Range("B6").Value = 4.227778939859434
$ws.Range("C6").Value = 0.1239377262776742
$ws.Range("D6").Value = 0.4616337648811282
$ws.Range("E6").Value = 0.118945123351784
$ws.Range("G6").Value = 0.002796559923334385
$ws.Range("J6").Value = 0.02569008908443049
$ws.Range("K6").Value = 3.514106071443564
$ws.Range("L6").Value = 0.6885660829091194
$ws.Range("N6").Value = 7.249103836207325

$ws.Range("B7").Value = 4.238026903427169
$ws.Range("C7").Value = 0.1255431557633671
$ws.Range("D7").Value = 0.4618441494702665
$ws.Range("E7").Value = 0.118848751742421
$ws.Range("G7").Value = 0.002794643017527341
$ws.Range("J7").Value = 0.02569971268182858
$ws.Range("K7").Value = 3.525538237139244
$ws.Range("L7").Value = 0.6891397863217605
$ws.Range("N7").Value = 7.270222644458499

$ws.Range("B8").Value = 4.28802251304927
$ws.Range("C8").Value = 0.1328208912899811
$ws.Range("D8").Value = 0.4632081520917666
$ws.Range("E8").Value = 0.1185090261120436
$ws.Range("G8").Value = 0.002786618768801347
$ws.Range("J8").Value = 0.02574196616005331
$ws.Range("K8").Value = 3.580175194879416
$ws.Range("L8").Value = 0.6923574591022685
$ws.Range("N8").Value = 7.364431963321096

$ws.Range("B9").Value = 4.403279503734893
$ws.Range("C9").Value = 0.1477993547934489
$ws.Range("D9").Value = 0.4674496605094731
$ws.Range("E9").Value = 0.1181495686865404
$ws.Range("G9").Value = 0.002772433632988904
$ws.Range("J9").Value = 0.025824156969664
$ws.Range("K9").Value = 3.702445996263862
$ws.Range("L9").Value = 0.7011371109497446
$ws.Range("N9").Value = 7.552909579680914

$ws.Range("B10").Value = 4.498530927011984
$ws.Range("C10").Value = 0.1592479425972328
$ws.Range("D10").Value = 0.4715233940346053
$ws.Range("E10").Value = 0.1180719960336916
$ws.Range("G10").Value = 0.002762945809134039
$ws.Range("J10").Value = 0.02588423441376086
$ws.Range("K10").Value = 3.80158810044577
$ws.Range("L10").Value = 0.7090994239574968
$ws.Range("N10").Value = 7.693753272535332

$ws.Range("B11").Value = 4.54417715940474
$ws.Range("C11").Value = 0.1645555236559346
$ws.Range("D11").Value = 0.4735851852650654
$ws.Range("E11").Value = 0.1180771209339468
$ws.Range("G11").Value = 0.002758829974264903
$ws.Range("J11").Value = 0.02591152027793164
$ws.Range("K11").Value = 3.848732362946976
$ws.Range("L11").Value = 0.7130515255789049
$ws.Range("N11").Value = 7.758362243348472

$ws.Range("B12").Value = 4.561796459475886
$ws.Range("C12").Value = 0.166579898838819
$ws.Range("D12").Value = 0.4743959700213338
$ws.Range("E12").Value = 0.1180848647622348
$ws.Range("G12").Value = 0.002757300023992881
$ws.Range("J12").Value = 0.02592184815526011
$ws.Range("K12").Value = 3.866880014880167
$ws.Range("L12").Value = 0.7145956390882731
$ws.Range("N12").Value = 7.782906815078036

$ws.Range("B13").Value = 4.557986953398313
$ws.Range("C13").Value = 0.1661432650057293
$ws.Range("D13").Value = 0.4742200169780233
$ws.Range("E13").Value = 0.1180829390205176
$ws.Range("G13").Value = 0.002757628255482948
$ws.Range("J13").Value = 0.02591962404464532
$ws.Range("K13").Value = 3.862958446990547
$ws.Range("L13").Value = 0.7142609709418508
$ws.Range("N13").Value = 7.777617190483681

$ws.Range("B14").Value = 4.545620008601929
$ws.Range("C14").Value = 0.164721778314231
$ws.Range("D14").Value = 0.4736512870119611
$ws.Range("E14").Value = 0.1180776417652929
$ws.Range("G14").Value = 0.002758703531634547
$ws.Range("J14").Value = 0.02591237004318181
$ws.Range("K14").Value = 3.850219458987056
$ws.Range("L14").Value = 0.7131776074870402
$ws.Range("N14").Value = 7.760379959154648

$ws.Range("B15").Value = 4.538088434248152
$ws.Range("C15").Value = 0.1638529720054294
$ws.Range("D15").Value = 0.4733068350660687
$ws.Range("E15").Value = 0.1180751525470072
$ws.Range("G15").Value = 0.002759365892478844
$ws.Range("J15").Value = 0.02590792619514026
$ws.Range("K15").Value = 3.842454934069451
$ws.Range("L15").Value = 0.7125202095722614
$ws.Range("N15").Value = 7.749831921265411

$ws.Range("B16").Value = 4.495594527921469
$ws.Range("C16").Value = 0.1589030995037035
$ws.Range("D16").Value = 0.471392851656887
$ws.Range("E16").Value = 0.1180724736126635
$ws.Range("G16").Value = 0.002763218803583027
$ws.Range("J16").Value = 0.02588245047524396
$ws.Range("K16").Value = 3.798548345910888
$ws.Range("L16").Value = 0.7088477924420573
$ws.Range("N16").Value = 7.689541855627112

$ws.Range("B17").Value = 4.470119773587271
$ws.Range("C17").Value = 0.1558921489825877
$ws.Range("D17").Value = 0.4702721391974478
$ws.Range("E17").Value = 0.1180811749614339
$ws.Range("G17").Value = 0.00276563360407164
$ws.Range("J17").Value = 0.02586681188735795
$ws.Range("K17").Value = 3.772137451058029
$ws.Range("L17").Value = 0.7066794639809615
$ws.Range("N17").Value = 7.652694564617605

$ws.Range("B18").Value = 4.455685285686229
$ws.Range("C18").Value = 0.1541696892992945
$ws.Range("D18").Value = 0.4696471703413465
$ws.Range("E18").Value = 0.1180899840849996
$ws.Range("G18").Value = 0.002767041389508302
$ws.Range("J18").Value = 0.02585781277802823
$ws.Range("K18").Value = 3.757138920834734
$ws.Range("L18").Value = 0.7054633551519629
$ws.Range("N18").Value = 7.631551592561095

$ws.Range("B19").Value = 4.450835412480046
$ws.Range("C19").Value = 0.1535880957226539
$ws.Range("D19").Value = 0.4694389381863715
$ws.Range("E19").Value = 0.1180936203581471
$ws.Range("G19").Value = 0.002767521285008821
$ws.Range("J19").Value = 0.02585476506389917
$ws.Range("K19").Value = 3.752093670621719
$ws.Range("L19").Value = 0.70505693307733
$ws.Range("N19").Value = 7.624401604597381

$ws.Range("B20").Value = 4.472809043686425
$ws.Range("C20").Value = 0.1562116997574492
$ws.Range("D20").Value = 0.4703894085750591
$ws.Range("E20").Value = 0.1180798550144093
$ws.Range("G20").Value = 0.002765374594126792
$ws.Range("J20").Value = 0.02586847706539608
$ws.Range("K20").Value = 3.774929023055961
$ws.Range("L20").Value = 0.7069070714409804
$ws.Range("N20").Value = 7.65661177506297

$ws.Range("B21").Value = 4.549243405197899
$ws.Range("C21").Value = 0.1651389080768979
$ws.Range("D21").Value = 0.4738175215645839
$ws.Range("E21").Value = 0.1180790402654424
$ws.Range("G21").Value = 0.002758386921357332
$ws.Range("J21").Value = 0.02591450083361657
$ws.Range("K21").Value = 3.853953188583546
$ws.Range("L21").Value = 0.7134945267216466
$ws.Range("N21").Value = 7.765440810787481

$ws.Range("B22").Value = 4.601145364280171
$ws.Range("C22").Value = 0.1710580188171491
$ws.Range("D22").Value = 0.4762330398400536
$ws.Range("E22").Value = 0.1181123294011108
$ws.Range("G22").Value = 0.002753986857315522
$ws.Range("J22").Value = 0.0259445534631908
$ws.Range("K22").Value = 3.907320960454854
$ws.Range("L22").Value = 0.7180769423392235
$ws.Range("N22").Value = 7.837025296995535

$ws.Range("B23").Value = 4.57326574350202
$ws.Range("C23").Value = 0.167891066650526
$ws.Range("D23").Value = 0.4749278051682921
$ws.Range("E23").Value = 0.1180914702887996
$ws.Range("G23").Value = 0.002756320047907387
$ws.Range("J23").Value = 0.02592851570728882
$ws.Range("K23").Value = 3.878679700667192
$ws.Range("L23").Value = 0.7156058341279419
$ws.Range("N23").Value = 7.798776988659881

$ws.Range("B24").Value = 4.471592566638378
$ws.Range("C24").Value = 0.1560672041542546
$ws.Range("D24").Value = 0.4703363308399986
$ws.Range("E24").Value = 0.118080439905091
$ws.Range("G24").Value = 0.002765491631816859
$ws.Range("J24").Value = 0.02586772426438344
$ws.Range("K24").Value = 3.773666375826963
$ws.Range("L24").Value = 0.706804075138777
$ws.Range("N24").Value = 7.654840676517296

$ws.Range("B25").Value = 4.37024856146337
$ws.Range("C25").Value = 0.1436702793031088
$ws.Range("D25").Value = 0.4661342564425723
$ws.Range("E25").Value = 0.1182140233331754
$ws.Range("G25").Value = 0.002776106264694408
$ws.Range("J25").Value = 0.02580198353903107
$ws.Range("K25").Value = 3.667740066449483
$ws.Range("L25").Value = 0.6984969450156768
$ws.Range("N25").Value = 7.501512109781146
